$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) header text of H1:J1 before anything is
# overwritten.
$oldH = $ws.Range("H1").Value()
$oldI = $ws.Range("I1").Value()
$oldJ = $ws.Range("J1").Value()

# Shift the last three header columns (گروه, توضیحات, شناسه) two columns
# to the right, from H:J to J:L, which opens up H:I for the two brand
# new "sale" fields. Writing the shifted-to cells first (while H1/I1
# still hold the original values) lets the shared-string table reuse the
# existing <si> entries for J1/K1/L1 instead of minting duplicates.
$ws.Range("J1").Value = $oldH
$ws.Range("K1").Value = $oldI
$ws.Range("L1").Value = $oldJ

# Now write the two new "sale" columns into the vacated H1/I1 cells.
$ws.Range("H1").Value = "فروش 2"
$ws.Range("I1").Value = "فروش 3"

# Match the number format Excel auto-applied to the two new header
# cells (integer format, style index 2 in styles.xml) -- same format
# already used by the neighbouring purchase/sale quantity columns.
$ws.Range("H1:I1").NumberFormat = "0"

# The column that used to hold the wide "توضیحات" (description) text
# keeps column index 9, but since that slot is now a short "sale"
# figure column, narrow it down.
$ws.Columns.Item(9).ColumnWidth = 10.3

# Reflect the new last-edited cell.
$ws.Range("J2").Select()
